$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 174, shifting rows 174:241 down to 175:242
$ws.Rows.Item(174).Insert()

# Populate the new row 174 with the new entry
$ws.Cells.Item(174, 1).Value = "basal diameter"
$ws.Cells.Item(174, 2).Value = "QUALITY"
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(174, 4).Value = 1
